# Updating filtered feeds from workflow
# Adds a new row (88) to the "Filtered Feeds" sheet with a new link/keywords/title,
# mirroring the formatting (hyperlink style) used by the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 88

$link = "https://www.360dx.com/clinical-lab-management/legal-ruling-fda-authority-over-ldts-raises-questions-advanced-dx-cdx"
$keywords = "CDx"
$title = "Legal Ruling on FDA Authority Over LDTs Raises Questions for Advanced Dx, CDx Developers"

# Column A: link (as a real hyperlink, matching the style already used in column A)
$ws.Cells.Item($newRow, 1).Value = $link
$ws.Hyperlinks.Add($ws.Cells.Item($newRow, 1), $link)
$ws.Cells.Item($newRow, 1).Style = $ws.Cells.Item($newRow - 1, 1).Style

# Column B: keywords
$ws.Cells.Item($newRow, 2).Value = $keywords

# Column C: title
$ws.Cells.Item($newRow, 3).Value = $title

Write-Host "Added row $newRow to $($ws.Name)"
